# Duplicate the last status-update sheet ("Status Update 25") to create a new
# "Status Update 26" sheet (lesson07), then update its narrative cells with
# this week's comments so they can be referenced later once coding begins.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Status Update 25")

# Copy the template sheet, placing the new copy right after the source sheet
# (i.e. at the very end of the workbook). This preserves all styles, merged
# cells, column widths, page setup, etc.
[void]$src.Copy($null, $src)

$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Status Update 26"

# --- Update the hours-tracking numbers for this week ---
$new.Range("J7").Value = 76
$new.Range("K7").Value = 76

# --- Update the narrative / status-comment cells with this week's text ---
$new.Range("C11").Value = "Achieved movement, the ability to change the ambient lighting color, and the ability to change the direction the lighting is coming from."
$new.Range("C12").Value = "Still have quite a bit of design to go"
$new.Range("E13").Value = "Going to have most of my research/design if not all finished 2 weeks from now, and start on coding the actual project."
$new.Range("C14").Value = "Actual hours since last update: 4 hours"
$new.Range("C16").Value = "Things are progressing well so far."
$new.Range("C17").Value = "Alma 43:7 - ...that he might bring [webgl] into subjection to the accomplishment of his designs."

# The new sheet becomes the active tab, with the cursor left on the last
# comment cell that was just filled in.
[void]$new.Activate()
[void]$new.Range("C17").Select()

# The previously-active sheet is no longer selected; leave its cursor on the
# first comment cell.
[void]$src.Activate()
[void]$src.Range("C11").Select()

# Re-activate the new sheet so it is the one shown/active when the file is
# saved (matches the workbook-level activeTab pointing at the new sheet).
[void]$new.Activate()

$excel.ActiveWindow.TabRatio = 987
